$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newStrings = @(
  "SCRIPT/P19P01A/enter07.ssb",
  "SCRIPT/D54P11A/enter06.ssb",
  "SCRIPT/D55P11A/enter06.ssb",
  "SCRIPT/D56P11A/enter05.ssb",
  "SCRIPT/D46P11A/enter02.ssb",
  "SCRIPT/D47P11A/enter02.ssb",
  "SCRIPT/D48P11A/enter02.ssb",
  "SCRIPT/D50P11A/enter02.ssb",
  "SCRIPT/D51P11A/enter03.ssb",
  "SCRIPT/D52P11A/enter02.ssb",
  "SCRIPT/D53P11A/enter02.ssb"
)

$startRow = 23
for ($i = 0; $i -lt $newStrings.Length; $i++) {
  $row = $startRow + $i
  $ws.Cells.Item($row, 1).Value = $newStrings[$i]
}

$ws.Range("C9").Select()
